$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell while keeping it a TEXT
# string (matching the sheet's existing inline-string "Price" column), by
# routing through a text formula and then flattening it to a static value
# via copy / paste-special-values. A direct .Value assignment of a
# numeric-looking string gets auto-coerced by Excel into a real number,
# which would change the cell's type/formatting from the original.
function Set-TextValue($range, $val) {
    $escaped = $val -replace '"', '""'
    $range.Formula = "=""" + $escaped + """"
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# --- Simple price refreshes ---
Set-TextValue $ws.Range("D2")  "246.36"
Set-TextValue $ws.Range("D3")  "22.26"
Set-TextValue $ws.Range("D4")  "5.352"
Set-TextValue $ws.Range("D5")  "0.05860"
Set-TextValue $ws.Range("D7")  "6.387"
Set-TextValue $ws.Range("D8")  "0.8134"
Set-TextValue $ws.Range("D9")  "0.9972"
Set-TextValue $ws.Range("D10") "0.1419"
Set-TextValue $ws.Range("D11") "0.03724"
Set-TextValue $ws.Range("D12") "0.07326"
Set-TextValue $ws.Range("D13") "0.03003"
Set-TextValue $ws.Range("D14") "4.175"
Set-TextValue $ws.Range("D15") "0.09398"
Set-TextValue $ws.Range("D16") "0.001600"
Set-TextValue $ws.Range("D17") "0.04817"
Set-TextValue $ws.Range("D18") "0.0005889"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue $ws.Range("D19") "0.005997"
Set-TextValue $ws.Range("D20") "0.004083"
Set-TextValue $ws.Range("D21") "0.0009882"
Set-TextValue $ws.Range("D23") "3.690"
Set-TextValue $ws.Range("D24") "2.232"
Set-TextValue $ws.Range("D25") "0.3247"
Set-TextValue $ws.Range("D27") "0.0002471"
Set-TextValue $ws.Range("D40") "0.03855"

# --- Rows 41-43: the three coins rotated positions (Kick -> row41,
#     BKEX -> row42, CEJI -> row43), each also getting an updated price ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006391"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.002600"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue $ws.Range("D44") "0.005215"
Set-TextValue $ws.Range("D45") "0.00005657"
Set-TextValue $ws.Range("D47") "0.7219"
Set-TextValue $ws.Range("D48") "0.08474"
